$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B47").Value = "Headers for Xbee"
$ws.Range("C47").Value = "NPPN101BFCN-RC"
$ws.Range("D47").Value = 10
$ws.Range("E47").Value = 1.018
$ws.Range("E47").NumberFormat = """$""#,##0.00"
$ws.Range("G47").Value = "http://www.digikey.com/product-detail/en/NPPN101BFCN-RC/S5751-10-ND/804812"

$ws.Range("C47").Font.Name = "Arial"
$ws.Range("C47").Font.Size = 12
$ws.Range("C47").Font.Color = 0

$ws.Rows(47).RowHeight = 15
